$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (target cluster: ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.022418
$ws.Range("H2").Value = 0.06725399999999999
$ws.Range("M2").Value = 0.3284223333333333
$ws.Range("N2").Value = 0.985267
$ws.Range("O2").Value = 0.04541528350839906
$ws.Range("P2").Value = 0.04541528350839906
$ws.Range("Q2").Value = 0.007362571868666666
$ws.Range("R2").Value = 0.066263146818
$ws.Range("S2").Value = 0.04541528350839906
$ws.Range("T2").Value = 0.04541528350839906

# Row 3 (target cluster: FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.022418
$ws.Range("H3").Value = 0.06725399999999999
$ws.Range("M3").Value = 4.062688333333334
$ws.Range("O3").Value = 0.5618014481290817
$ws.Range("P3").Value = 0.5618014481290816
$ws.Range("Q3").Value = 0.09107734705666666
$ws.Range("R3").Value = 0.8196961235100001
$ws.Range("S3").Value = 0.5618014481290817
$ws.Range("T3").Value = 0.5618014481290816

# Row 4 (target cluster: MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.022418
$ws.Range("H4").Value = 0.06725399999999999
$ws.Range("M4").Value = 2.840427
$ws.Range("N4").Value = 8.521281
$ws.Range("O4").Value = 0.3927832683625193
$ws.Range("P4").Value = 0.3927832683625193
$ws.Range("Q4").Value = 0.063676692486
$ws.Range("R4").Value = 0.5730902323739999
$ws.Range("S4").Value = 0.3927832683625193
$ws.Range("T4").Value = 0.3927832683625193
